# V1.0, final version submited to OC
# Update the progress percentages on the Dashboard sheet to 100% (1) for the
# remaining tasks, then move the active selection to D14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("C21").Value = 1

$ws.Calculate() | Out-Null

$ws.Range("D14").Select() | Out-Null
